# Updated cryptos list with latest prices and 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.376.45"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.451.71"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.45"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.20"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "2.449.20"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  +2.26%  "
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.346"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.32"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "2.897.33"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "62.187.91"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "2.448.86"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.85"
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.12"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.94"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("E23").Value = "  -6.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.64"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.27"
$ws.Range("E26").Value = "  +3.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "588.56"
$ws.Range("E27").Value = "  -4.66%  "
$ws.Range("D29").Value = "0.0₃0960"
$ws.Range("E29").Value = "  -2.96%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -3.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.03"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.137"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.92"
$ws.Range("E35").Value = "  -2.93%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.45"
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.379"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.87"
$ws.Range("E39").Value = "  +4.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.33"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.39"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.04"
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.73"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.45"
$ws.Range("E45").Value = "  -4.00%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.38"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.65"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0251"
$ws.Range("E48").Value = "  +13.25%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.607"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.77"
$ws.Range("E51").Value = "  -3.55%  "
